# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C2:C14) from 2023-10-13 (45212)
# to 2023-10-22 (45221) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C14").Value = 45221
